# Update cryptos list price/volume data (GitHub Actions scheduled refresh).
# Also swaps rank of Uniswap and Wrapped liquid staked Ether 2.0 (rows 20/21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column values are free-form text ("30.294.89", "1.001", ...), not
# real numbers. A leading apostrophe keeps values that otherwise parse as
# numbers (e.g. "1.001") stored as literal text, same as typing them into
# Excel by hand; plain text values (links, names, percents) need no prefix.
$ws.Range("D2").Value = "30.294.89"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.882.64"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'235.68"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4683"
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("D8").Value = "'0.2830"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "'0.06590"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").Value = "'20.60"
$ws.Range("E10").Value = "  +9.55%  "
$ws.Range("D11").Value = "'0.07768"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "'97.64"
$ws.Range("E12").Value = "  -3.99%  "
$ws.Range("D13").Value = "1.885.97"
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("D14").Value = "'5.070"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").Value = "'0.6731"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "'283.59"
$ws.Range("E16").Value = "  +6.91%  "
$ws.Range("D17").Value = "30.311.49"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'12.61"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'5.396"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.114.63"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "'0.000007274"
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'6.166"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "'9.369"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "'168.02"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Value = "'19.19"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "'1.987"
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "'0.09680"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("D31").Value = "'4.371"
$ws.Range("E31").Value = "  -7.52%  "
$ws.Range("D32").Value = "'1.476"
$ws.Range("E32").Value = "  -2.15%  "
$ws.Range("D33").Value = "'4.118"
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("D34").Value = "'0.04669"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").Value = "'0.7050"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("D36").Value = "'1.096"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").Value = "'2.713"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'0.01868"
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("D39").Value = "'6.569"
$ws.Range("E39").Value = "  +5.41%  "
$ws.Range("D40").Value = "'2.523"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("D41").Value = "'71.89"
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").Value = "'1.963"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "'0.8648"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D45").Value = "'102.92"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("D46").Value = "'0.4178"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("D47").Value = "'986.04"
$ws.Range("E47").Value = "  +7.50%  "
$ws.Range("D48").Value = "'7.262"
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("D49").Value = "'9.185"
$ws.Range("E49").Value = "  +4.88%  "
$ws.Range("D50").Value = "'33.86"
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("D51").Value = "'0.1144"
$ws.Range("E51").Value = "  -4.86%  "
